$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '25.787.26'
$ws.Range("E2").Value = '  -2.58%  '

# Row 3
$ws.Range("D3").Value = '1.745.12'
$ws.Range("E3").Value = '  -5.00%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.30'
$ws.Range("E5").Value = '  -9.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5050'
$ws.Range("E7").Value = '  -6.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.92'
$ws.Range("E8").Value = '  -6.57%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2654'
$ws.Range("E9").Value = '  -12.01%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06151'
$ws.Range("E10").Value = '  -10.42%  '

# Row 11
$ws.Range("D11").Value = '1.745.30'
$ws.Range("E11").Value = '  -5.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06919'
$ws.Range("E12").Value = '  -4.29%  '

# Row 13
$ws.Range("E13").Value = '  -12.41%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.504'
$ws.Range("E14").Value = '  -9.61%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5992'
$ws.Range("E15").Value = '  -18.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.02'
$ws.Range("E16").Value = '  -13.84%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.08%  '

# Row 19
$ws.Range("D19").Value = '25.793.17'
$ws.Range("E19").Value = '  -2.66%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006861'
$ws.Range("E20").Value = '  -12.82%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.59'
$ws.Range("E21").Value = '  -16.05%  '

# Row 22
$ws.Range("D22").Value = '1.967.13'
$ws.Range("E22").Value = '  -5.49%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.049'
$ws.Range("E23").Value = '  -11.74%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.209'
$ws.Range("E24").Value = '  -12.65%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.140'
$ws.Range("E25").Value = '  -11.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '137.45'
$ws.Range("E26").Value = '  -3.73%  '

# Row 27
$ws.Range("E27").Value = '  -9.93%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.818'
$ws.Range("E28").Value = '  -17.16%  '

# Row 29
$ws.Range("E29").Value = '  -11.46%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '103.64'
$ws.Range("E30").Value = '  -6.33%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.770'
$ws.Range("E31").Value = '  -10.72%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08108'
$ws.Range("E32").Value = '  -8.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.465'
$ws.Range("E33").Value = '  -13.72%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04510'
$ws.Range("E34").Value = '  -6.16%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9994'
$ws.Range("E35").Value = '  -0.07%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.646'
$ws.Range("E36").Value = '  -9.27%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9841'
$ws.Range("E37").Value = '  -12.85%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6071'
$ws.Range("E38").Value = '  -16.54%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.663'
$ws.Range("E39").Value = '  -13.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01549'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.920'
$ws.Range("E41").Value = '  -15.29%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9999'
$ws.Range("E42").Value = '  -0.04%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.37'
$ws.Range("E43").Value = '  -4.96%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3817'
$ws.Range("E44").Value = '  -18.91%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.080'
$ws.Range("E45").Value = '  -13.53%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7332'
$ws.Range("E46").Value = '  -18.91%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05360'
$ws.Range("E47").Value = '  -7.16%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1109'
$ws.Range("E48").Value = '  -10.51%  '

# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.929'
$ws.Range("E49").Value = '  -19.44%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.17'
$ws.Range("E50").Value = '  -13.14%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.64'
$ws.Range("E51").Value = '  -12.12%  '
